# Updates the "cryptos" price/volume table with freshly scraped values.
# Note: several Price-column values look like plain numbers (e.g. "1.00",
# "0.000329"). The source cells are plain text, so a leading apostrophe
# (PowerShell '' inside a single-quoted string = one literal ') is used
# to force Excel to keep them as text instead of silently converting them
# to numeric values (which would drop formatting like trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.926.97'
$ws.Range('E2').Value = '  +2.96%  '
$ws.Range('D3').Value = '3.801.12'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''700.30'
$ws.Range('E5').Value = '  +8.93%  '
$ws.Range('D6').Value = '''172.74'
$ws.Range('E6').Value = '  +4.40%  '
$ws.Range('D7').Value = '3.801.17'
$ws.Range('E7').Value = '  +0.97%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('E10').Value = '  +2.99%  '
$ws.Range('D11').Value = '''7.41'
$ws.Range('E11').Value = '  +6.76%  '
$ws.Range('D12').Value = '''0.461'
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').Value = '''0.0000259'
$ws.Range('E13').Value = '  +8.40%  '
$ws.Range('D14').Value = '''36.45'
$ws.Range('E14').Value = '  +4.59%  '
$ws.Range('D15').Value = '4.440.42'
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('D16').Value = '3.805.65'
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('D17').Value = '70.860.71'
$ws.Range('E17').Value = '  +2.85%  '
$ws.Range('D18').Value = '''17.89'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('E19').Value = '  +3.08%  '
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').Value = '''11.08'
$ws.Range('E21').Value = '  +16.04%  '
$ws.Range('D22').Value = '''482.91'
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('E23').Value = '  +1.46%  '
$ws.Range('D24').Value = '''84.43'
$ws.Range('E24').Value = '  +3.27%  '
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('D26').Value = '''12.45'
$ws.Range('E26').Value = '  +2.51%  '
$ws.Range('E27').Value = '  +3.90%  '
$ws.Range('D28').Value = '''10.48'
$ws.Range('E28').Value = '  +3.97%  '
$ws.Range('D29').Value = '3.950.42'
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('D31').Value = '''3.11'
$ws.Range('E31').Value = '  +16.06%  '
$ws.Range('D32').Value = '''7.55'
$ws.Range('E32').Value = '  +6.07%  '
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('D34').Value = '''0.185'
$ws.Range('E34').Value = '  +7.18%  '
$ws.Range('E35').Value = '  +3.63%  '
$ws.Range('E36').Value = '  +4.06%  '
$ws.Range('D37').Value = '''1.00'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  +2.30%  '
$ws.Range('D39').Value = '''3.45'
$ws.Range('E39').Value = '  +6.64%  '
$ws.Range('E40').Value = '  +4.93%  '
$ws.Range('E41').Value = '  +12.23%  '
$ws.Range('B42').Value = 'FLOKI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D42').Value = '''0.000329'
$ws.Range('E42').Value = '  +23.58%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '''0.978'
$ws.Range('E43').Value = '  +2.21%  '
$ws.Range('D44').Value = '''0.999'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = '''162.16'
$ws.Range('E46').Value = '  +4.27%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = '''45.64'
$ws.Range('E47').Value = '  +1.78%  '
$ws.Range('D48').Value = '''48.86'
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('D49').Value = '''0.303'
$ws.Range('E49').Value = '  +2.49%  '
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('E51').Value = '  +2.68%  '
